# Weekly refresh of "Haba" price data for Comercializadora del Agro de Limarí
# - Updates Fecha (D) and Volumen/Precio columns (J,K,L,M,P) for existing rows 2-34
# - Inserts 3 new data rows above the former row 35 (new rows 35-37), pushing the
#   former row 35 down to row 38 intact

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{ Row=2; D=44329; J=1000; K=12000; L=13000; M=12500; P=500 },
  @{ Row=3; D=44406; J=800; K=10000; L=11000; M=10500; P=420 },
  @{ Row=4; D=44455; J=600 },
  @{ Row=6; D=44441; J=1100; K=11000; L=12000; M=11500; P=460 },
  @{ Row=7; D=44420; J=1000; K=10000; L=11000; M=10500; P=420 },
  @{ Row=8; D=44462; J=800; K=9000; L=10000; M=9500; P=380 },
  @{ Row=9; D=44371; J=500; L=12000; M=11000; P=440 },
  @{ Row=10; D=44399; J=500; K=9000; L=10000; M=9500; P=380 },
  @{ Row=12; D=44419; J=1100; K=11000; L=12000; M=11500; P=460 },
  @{ Row=13; D=44412; J=1000 },
  @{ Row=14; D=44336; J=1200; K=12000; L=13000; M=12500; P=500 },
  @{ Row=15; D=44308; J=400; K=11000; L=12000; M=11500; P=460 },
  @{ Row=16; D=44343 },
  @{ Row=17; D=44447; J=1000; K=10000; M=11000; P=440 },
  @{ Row=18; D=44356; J=1000; K=11000; L=12000; M=11500; P=460 },
  @{ Row=19; D=44454; J=800 },
  @{ Row=20; D=44427; J=360; K=10000; L=11000; M=10500; P=420 },
  @{ Row=21; D=44413; K=10000; L=11000; M=10500; P=420 },
  @{ Row=23; D=44469; J=600; K=5000; L=6000; M=5500; P=220 },
  @{ Row=24; D=44434; K=10000; L=11000; M=10500; P=420 },
  @{ Row=25; D=44391; J=500; K=9000; L=10000; M=9500; P=380 },
  @{ Row=26; D=44335; K=12000; L=13000; M=12500; P=500 },
  @{ Row=27; D=44475; J=1200; K=5000; L=6000; M=5500; P=220 },
  @{ Row=28; D=44349; J=600; L=12000; M=11000; P=440 },
  @{ Row=29; D=44448; J=800; K=10000; L=12000; M=11000; P=440 },
  @{ Row=30; D=44385; J=600; K=8000; L=9000; M=8500; P=340 },
  @{ Row=31; D=44435; J=600; K=10000; L=11000; M=10500; P=420 },
  @{ Row=32; D=44377; J=800; K=9000; L=10000; M=9500; P=380 },
  @{ Row=33; D=44468; J=700; K=5000; L=6000; M=5500; P=220 },
  @{ Row=34; D=44384; J=700; K=8000; L=9000; M=8500; P=340 },
)

foreach ($u in $updates) {
  $r = $u.Row
  if ($u.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $u.D }
  if ($u.ContainsKey("J")) { $ws.Cells.Item($r, 10).Value = $u.J }
  if ($u.ContainsKey("K")) { $ws.Cells.Item($r, 11).Value = $u.K }
  if ($u.ContainsKey("L")) { $ws.Cells.Item($r, 12).Value = $u.L }
  if ($u.ContainsKey("M")) { $ws.Cells.Item($r, 13).Value = $u.M }
  if ($u.ContainsKey("P")) { $ws.Cells.Item($r, 16).Value = $u.P }
}

# --- Insert 3 new rows above row 35 (pushes former row 35 down to row 38) ---
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()

$newRows = @(
  @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44363; E=4; F=100112026; G="Haba"; H="Sin especificar"; I="Primera"; J=900; K=11000; L=12000; M=11500; N="`$/saco 25 kilos"; O="Provincia de Limarí"; P=460; Q=25; R="Hortaliza" },
  @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44328; E=4; F=100112026; G="Haba"; H="Sin especificar"; I="Primera"; J=900; K=11000; L=12000; M=11500; N="`$/saco 25 kilos"; O="Provincia de Limarí"; P=460; Q=25; R="Hortaliza" },
  @{ A=2; B="Comercializadora del Agro de Limarí"; C="Coquimbo"; D=44392; E=4; F=100112026; G="Haba"; H="Sin especificar"; I="Primera"; J=600; K=9000; L=10000; M=9500; N="`$/saco 25 kilos"; O="Provincia de Limarí"; P=380; Q=25; R="Hortaliza" },
)

$newRowNums = @(35, 36, 37)
for ($i = 0; $i -lt $newRowNums.Count; $i++) {
  $r = $newRowNums[$i]
  $u = $newRows[$i]
  for ($c = 1; $c -le 18; $c++) {
    $letter = [char](64 + $c)
    $ws.Cells.Item($r, $c).Value = $u[[string]$letter]
  }
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()